$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1092.4445
$ws.Range("I15").Value = 1092.4445
$ws.Range("K15").Value = 3277.3335
$ws.Range("M15").Value = -3108.3335

$ws.Range("H43").Value = 1898.3334
$ws.Range("J43").Value = 1898.3334
$ws.Range("L43").Value = 1898.3334
$ws.Range("N43").Value = -2036.3334

$ws.Range("H57").Value = 49983
$ws.Range("J57").Value = 49983
$ws.Range("L57").Value = 149949
$ws.Range("N57").Value = -150947

$ws.Range("H80").Value = 6716.3
$ws.Range("I80").Value = 4037
$ws.Range("K80").Value = 12111
$ws.Range("M80").Value = -11113

$ws.Range("H83").Value = 6716.3
$ws.Range("I83").Value = 4037
$ws.Range("K83").Value = 36333
$ws.Range("M83").Value = -31341

$ws.Range("H99").Value = 1120.125
$ws.Range("I99").Value = 401.69232
$ws.Range("K99").Value = 1205.07696
$ws.Range("M99").Value = 292.9230400000001

$ws.Range("H111").Value = 1061.909
$ws.Range("I111").Value = 1327.1428
$ws.Range("K111").Value = 3981.4284
$ws.Range("M111").Value = -914.4284000000002

$ws.Range("H132").Value = 128428.875
$ws.Range("I132").Value = 128428.875
$ws.Range("K132").Value = 385286.625
$ws.Range("M132").Value = -382756.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 43100
$ws.Range("J56").Value = 43100
$ws.Range("L56").Value = 43100
$ws.Range("N56").Value = -44584

$ws.Range("H61").Value = 1761.75
$ws.Range("I61").Value = 1761.75
$ws.Range("K61").Value = 1761.75
$ws.Range("M61").Value = -1549.75

$ws.Range("H136").Value = 1761.75
$ws.Range("I136").Value = 1761.75
$ws.Range("K136").Value = 5285.25
$ws.Range("M136").Value = -2735.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2112.8
$ws.Range("I5").Value = 1721.6666
$ws.Range("J5").Value = 2699.5
$ws.Range("K5").Value = 1721.6666
$ws.Range("L5").Value = 2699.5
$ws.Range("M5").Value = -1608.6666
$ws.Range("N5").Value = -2925.5

$ws.Range("H22").Value = 589
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H86").Value = 2113.7778
$ws.Range("I86").Value = 2003.4286
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 2003.4286
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -880.4286
$ws.Range("N86").Value = -4746

$ws.Range("H89").Value = 2113.7778
$ws.Range("I89").Value = 2003.4286
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 10017.143
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -4401.143
$ws.Range("N89").Value = -23732

$ws.Range("H105").Value = 3719.875
$ws.Range("I105").Value = 3993.6
$ws.Range("J105").Value = 3263.6667
$ws.Range("K105").Value = 3993.6
$ws.Range("L105").Value = 3263.6667
$ws.Range("M105").Value = -2246.6
$ws.Range("N105").Value = -6757.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 33649.5
$ws.Range("I69").Value = 14000
$ws.Range("K69").Value = 14000
$ws.Range("M69").Value = -13251

$ws.Range("H72").Value = 33649.5
$ws.Range("I72").Value = 14000
$ws.Range("K72").Value = 42000
$ws.Range("M72").Value = -38256

$ws.Range("H132").Value = 2141.2856
$ws.Range("I132").Value = 2098.2
$ws.Range("J132").Value = 2249
$ws.Range("K132").Value = 6294.599999999999
$ws.Range("L132").Value = 6747
$ws.Range("M132").Value = -3764.599999999999
$ws.Range("N132").Value = -11807

$ws.Range("H134").Value = 2414.375
$ws.Range("I134").Value = 2414.375
$ws.Range("K134").Value = 7243.125
$ws.Range("M134").Value = -4708.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 146.73334
$ws.Range("I2").Value = 99.5
$ws.Range("K2").Value = 597
$ws.Range("M2").Value = -484

$ws.Range("H4").Value = 89166.734
$ws.Range("I4").Value = 947.03705
$ws.Range("J4").Value = 429442.72
$ws.Range("K4").Value = 2841.11115
$ws.Range("L4").Value = 1288328.16
$ws.Range("M4").Value = -2729.11115
$ws.Range("N4").Value = -1288552.16

$ws.Range("H69").Value = 2063.25
$ws.Range("J69").Value = 2251.4285
$ws.Range("L69").Value = 6754.2855
$ws.Range("N69").Value = -8376.2855

$ws.Range("H72").Value = 2063.25
$ws.Range("J72").Value = 2251.4285
$ws.Range("L72").Value = 20262.8565
$ws.Range("N72").Value = -28374.8565

$ws.Range("H80").Value = 4828.4443
$ws.Range("J80").Value = 4993.1665
$ws.Range("L80").Value = 14979.4995
$ws.Range("N80").Value = -16851.4995

$ws.Range("H81").Value = 2766
$ws.Range("I81").Value = 2140
$ws.Range("J81").Value = 4957
$ws.Range("K81").Value = 6420
$ws.Range("L81").Value = 14871
$ws.Range("M81").Value = -5297
$ws.Range("N81").Value = -17117

$ws.Range("H83").Value = 4828.4443
$ws.Range("J83").Value = 4993.1665
$ws.Range("L83").Value = 44938.4985
$ws.Range("N83").Value = -54298.4985

$ws.Range("H84").Value = 2766
$ws.Range("I84").Value = 2140
$ws.Range("J84").Value = 4957
$ws.Range("K84").Value = 19260
$ws.Range("L84").Value = 44613
$ws.Range("M84").Value = -13644
$ws.Range("N84").Value = -55845

$ws.Range("H94").Value = 3529.1177
$ws.Range("I94").Value = 2582.9167
$ws.Range("J94").Value = 5800
$ws.Range("K94").Value = 7748.750100000001
$ws.Range("L94").Value = 17400
$ws.Range("M94").Value = -7072.750100000001
$ws.Range("N94").Value = -18752

$ws.Range("H122").Value = 40023.848
$ws.Range("I122").Value = 962.5
$ws.Range("J122").Value = 73505
$ws.Range("K122").Value = 8662.5
$ws.Range("L122").Value = 661545
$ws.Range("M122").Value = -6212.5
$ws.Range("N122").Value = -666445

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 80749.875
$ws.Range("J18").Value = 24998
$ws.Range("L18").Value = 24998
$ws.Range("N18").Value = -25584

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws.Range("H132").Value = 4443.5264
$ws.Range("I132").Value = 4393.125
$ws.Range("J132").Value = 4712.3335
$ws.Range("K132").Value = 13179.375
$ws.Range("L132").Value = 14137.0005
$ws.Range("M132").Value = -10649.375
$ws.Range("N132").Value = -19197.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3418.3845
$ws.Range("I16").Value = 2267.9092
$ws.Range("K16").Value = 2267.9092
$ws.Range("M16").Value = -2097.9092

$ws.Range("H32").Value = 7671
$ws.Range("I32").Value = 7671
$ws.Range("K32").Value = 7671
$ws.Range("M32").Value = -7354

$ws.Range("H94").Value = 40330
$ws.Range("J94").Value = 40330
$ws.Range("L94").Value = 40330
$ws.Range("N94").Value = -41682

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9999.5
$ws.Range("I41").Value = 9999
$ws.Range("K41").Value = 9999
$ws.Range("M41").Value = -9609

$ws.Range("H126").Value = 2685.5
$ws.Range("I126").Value = 2049.75
$ws.Range("K126").Value = 6149.25
$ws.Range("M126").Value = -3679.25
